$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper: write a text-like value (e.g. "1.46%", "4.80") into a cell while
# keeping it stored as literal text (not auto-converted to a number) AND
# without leaving a stray NumberFormat/style behind on the cell.
function Set-TextValue {
    param(
        [string]$CellRef,
        [string]$Text
    )
    $rng = $ws.Range($CellRef)
    $rng.NumberFormat = "@"
    $rng.Value = $Text
    $rng.ClearFormats()
}

# --- Plain numeric data changes (xl/worksheets/sheet1.xml <v> values) ---

# Row 2: Cantidad de puntos
$ws.Range("D2").Value = 17040

# Row 4: Cierres del mes
$ws.Range("B4").Value = 204
$ws.Range("C4").Value = 159
$ws.Range("D4").Value = 363

# Row 7: En proceso de cierre
$ws.Range("B7").Value = 5
$ws.Range("C7").Value = 11
$ws.Range("D7").Value = 16

# Row 13: Puntos bloqueados - Activos
$ws.Range("B13").Value = 63
$ws.Range("C13").Value = 134
$ws.Range("D13").Value = 197

# Row 14: Puntos bloqueados - Inactivos
$ws.Range("B14").Value = 19
$ws.Range("C14").Value = 24
$ws.Range("D14").Value = 43

# --- Percentage/decimal text values (shared strings) ---

# Row 10: Puntos con malas practicas (%)  -- "ICX" related per commit message
Set-TextValue "B10" "1.46%"
Set-TextValue "C10" "1.12%"
# D10 unchanged (1.95%)

# Row 15: Puntos bloqueados - Activos (%)
Set-TextValue "B15" "26.25%"
Set-TextValue "C15" "55.83%"
Set-TextValue "D15" "82.08%"

# Row 16: Puntos bloqueados - Inactivos (%)
Set-TextValue "B16" "7.92%"
Set-TextValue "C16" "10.0%"
Set-TextValue "D16" "17.92%"

# Row 17: ICX
Set-TextValue "B17" "4.80"
Set-TextValue "C17" "4.80"
Set-TextValue "D17" "4.80"

# Row 18: NPS
Set-TextValue "B18" "82.25%"
Set-TextValue "C18" "82.25%"
Set-TextValue "D18" "82.25%"
